# "Generate Report for Handoff"
#
# The localization-status report has moved from "In Translation" to
# "Ready for handoff" for both the zh-cn and de-de targets, and the
# timestamps that track when each handoff package was (re)generated have
# been bumped. Update all three sheets (Overview, zh-cn, de-de) to match,
# and widen the two "Status" columns that now hold the longer text.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row 2 holds the zh-cn status (E2), de-de status (F2)
# and the "Latest HO Xliff Generate Date" (G2).
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-08-28 14:58:30"

# Status columns got wider now that they hold "Ready for handoff"
# instead of "In Translation".
$overview.Columns("E:F").ColumnWidth = 16.3333333333333

# ---------------------------------------------------------------------
# zh-cn detail sheet: Status (C2) and Latest Handoff Datetime (H2).
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-08-28 14:58:25"
$zhcn.Columns("C").ColumnWidth = 16.3333333333333

# ---------------------------------------------------------------------
# de-de detail sheet: Status (C2) and Latest Handoff Datetime (H2).
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-08-28 14:58:30"
$dede.Columns("C").ColumnWidth = 16.3333333333333
